$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(145).Insert()

$ws.Cells.Item(145, 1).Value = 5
$ws.Cells.Item(145, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(145, 3).Value = "Maule"
$ws.Cells.Item(145, 4).Value = 44589
$ws.Cells.Item(145, 5).Value = 7
$ws.Cells.Item(145, 6).Value = 100114013
$ws.Cells.Item(145, 7).Value = "Zanahoria"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 300
$ws.Cells.Item(145, 11).Value = 8000
$ws.Cells.Item(145, 12).Value = 8000
$ws.Cells.Item(145, 13).Value = 8000
$ws.Cells.Item(145, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(145, 15).Value = "Región de Ñuble"
$ws.Cells.Item(145, 16).Value = 400
$ws.Cells.Item(145, 17).Value = 20
$ws.Cells.Item(145, 18).Value = "Hortaliza"
